$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 138, pushing existing rows 138-199 down to 139-200.
$ws.Rows.Item(138).Insert()

# Populate the new row 138 with the weekly record.
$ws.Range("A138").Value = 11
$ws.Range("B138").Value = "Vega Monumental Concepción"
$ws.Range("C138").Value = "Bíobío"
$ws.Range("D138").Value = 45029
$ws.Range("E138").Value = 8
$ws.Range("F138").Value = 100112043
$ws.Range("G138").Value = "Pepino ensalada"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 200
$ws.Range("K138").Value = 6500
$ws.Range("L138").Value = 7000
$ws.Range("M138").Value = 6750
$ws.Range("N138").Value = "`$/caja 60 unidades"
$ws.Range("O138").Value = "Región de Arica y Parinacota"
$ws.Range("P138").Value = 112
$ws.Range("Q138").Value = 60
$ws.Range("R138").Value = "Hortaliza"

# Make sure the date cell keeps the same number format as the rest of column D.
$ws.Range("D138").NumberFormat = $ws.Range("D139").NumberFormat
